$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object "object[,]" 1,20
$row2[0,0] = "ECs"
$row2[0,1] = "Clcf1"
$row2[0,2] = "Cntfr"
$row2[0,3] = "ECs"
$row2[0,4] = 3
$row2[0,5] = 1
$row2[0,6] = 1.874436
$row2[0,7] = 5.623308
$row2[0,8] = 0.1442186763702422
$row2[0,9] = 0.1442186763702422
$row2[0,10] = 1
$row2[0,11] = 0.3333333333333333
$row2[0,12] = 0.03844333333333334
$row2[0,13] = 0.11533
$row2[0,14] = 0.006478890266597937
$row2[0,15] = 0.006478890266597936
$row2[0,16] = 0.07205956796
$row2[0,17] = 0.64853611164
$row2[0,18] = 0.0009343769785967999
$row2[0,19] = 0.0009343769785967996
$ws.Range("A2:T2").Value = $row2

$row3 = New-Object "object[,]" 1,20
$row3[0,0] = "ECs"
$row3[0,1] = "Clcf1"
$row3[0,2] = "Cntfr"
$row3[0,3] = "FAPs"
$row3[0,4] = 3
$row3[0,5] = 1
$row3[0,6] = 1.874436
$row3[0,7] = 5.623308
$row3[0,8] = 0.1442186763702422
$row3[0,9] = 0.1442186763702422
$row3[0,10] = 3
$row3[0,11] = 1
$row3[0,12] = 5.666771333333333
$row3[0,13] = 17.000314
$row3[0,14] = 0.95502617622222
$row3[0,15] = 0.9550261762222199
$row3[0,16] = 10.622000190968
$row3[0,17] = 95.598001718712
$row3[0,18] = 0.1377326110337022
$row3[0,19] = 0.1377326110337022
$ws.Range("A3:T3").Value = $row3

$row4 = New-Object "object[,]" 1,20
$row4[0,0] = "ECs"
$row4[0,1] = "Clcf1"
$row4[0,2] = "Cntfr"
$row4[0,3] = "MuSCs"
$row4[0,4] = 3
$row4[0,5] = 1
$row4[0,6] = 1.874436
$row4[0,7] = 5.623308
$row4[0,8] = 0.1442186763702422
$row4[0,9] = 0.1442186763702422
$row4[0,10] = 3
$row4[0,11] = 1
$row4[0,12] = 0.2284146666666667
$row4[0,13] = 0.685244
$row4[0,14] = 0.03849493351118214
$row4[0,15] = 0.03849493351118213
$row4[0,16] = 0.428148674128
$row4[0,17] = 3.853338067151999
$row4[0,18] = 0.005551688357943168
$row4[0,19] = 0.005551688357943165
$ws.Range("A4:T4").Value = $row4

$row5 = New-Object "object[,]" 1,20
$row5[0,0] = "FAPs"
$row5[0,1] = "Clcf1"
$row5[0,2] = "Cntfr"
$row5[0,3] = "ECs"
$row5[0,4] = 3
$row5[0,5] = 1
$row5[0,6] = 3.197979
$row5[0,7] = 9.593937
$row5[0,8] = 0.2460517715407892
$row5[0,9] = 0.2460517715407892
$row5[0,10] = 1
$row5[0,11] = 0.3333333333333333
$row5[0,12] = 0.03844333333333334
$row5[0,13] = 0.11533
$row5[0,14] = 0.006478890266597937
$row5[0,15] = 0.006478890266597936
$row5[0,16] = 0.12294097269
$row5[0,17] = 1.10646875421
$row5[0,18] = 0.001594142427714798
$row5[0,19] = 0.001594142427714798
$ws.Range("A5:T5").Value = $row5

$row6 = New-Object "object[,]" 1,20
$row6[0,0] = "FAPs"
$row6[0,1] = "Clcf1"
$row6[0,2] = "Cntfr"
$row6[0,3] = "FAPs"
$row6[0,4] = 3
$row6[0,5] = 1
$row6[0,6] = 3.197979
$row6[0,7] = 9.593937
$row6[0,8] = 0.2460517715407892
$row6[0,9] = 0.2460517715407892
$row6[0,10] = 3
$row6[0,11] = 1
$row6[0,12] = 5.666771333333333
$row6[0,13] = 17.000314
$row6[0,14] = 0.95502617622222
$row6[0,15] = 0.9550261762222199
$row6[0,16] = 18.122215721802
$row6[0,17] = 163.099941496218
$row6[0,18] = 0.2349858825273032
$row6[0,19] = 0.2349858825273031
$ws.Range("A6:T6").Value = $row6

$row7 = New-Object "object[,]" 1,20
$row7[0,0] = "FAPs"
$row7[0,1] = "Clcf1"
$row7[0,2] = "Cntfr"
$row7[0,3] = "MuSCs"
$row7[0,4] = 3
$row7[0,5] = 1
$row7[0,6] = 3.197979
$row7[0,7] = 9.593937
$row7[0,8] = 0.2460517715407892
$row7[0,9] = 0.2460517715407892
$row7[0,10] = 3
$row7[0,11] = 1
$row7[0,12] = 0.2284146666666667
$row7[0,13] = 0.685244
$row7[0,14] = 0.03849493351118214
$row7[0,15] = 0.03849493351118213
$row7[0,16] = 0.730465307292
$row7[0,17] = 6.574187765628
$row7[0,18] = 0.009471746585771258
$row7[0,19] = 0.009471746585771256
$ws.Range("A7:T7").Value = $row7

$row8 = New-Object "object[,]" 1,20
$row8[0,0] = "MuSCs"
$row8[0,1] = "Clcf1"
$row8[0,2] = "Cntfr"
$row8[0,3] = "ECs"
$row8[0,4] = 3
$row8[0,5] = 1
$row8[0,6] = 6.825289333333334
$row8[0,7] = 20.475868
$row8[0,8] = 0.5251361975000832
$row8[0,9] = 0.5251361975000832
$row8[0,10] = 1
$row8[0,11] = 0.3333333333333333
$row8[0,12] = 0.03844333333333334
$row8[0,13] = 0.11533
$row8[0,14] = 0.006478890266597937
$row8[0,15] = 0.006478890266597936
$row8[0,16] = 0.2623868729377778
$row8[0,17] = 2.36148185644
$row8[0,18] = 0.003402299798621541
$row8[0,19] = 0.003402299798621541
$ws.Range("A8:T8").Value = $row8

$row9 = New-Object "object[,]" 1,20
$row9[0,0] = "MuSCs"
$row9[0,1] = "Clcf1"
$row9[0,2] = "Cntfr"
$row9[0,3] = "FAPs"
$row9[0,4] = 3
$row9[0,5] = 1
$row9[0,6] = 6.825289333333334
$row9[0,7] = 20.475868
$row9[0,8] = 0.5251361975000832
$row9[0,9] = 0.5251361975000832
$row9[0,10] = 3
$row9[0,11] = 1
$row9[0,12] = 5.666771333333333
$row9[0,13] = 17.000314
$row9[0,14] = 0.95502617622222
$row9[0,15] = 0.9550261762222199
$row9[0,16] = 38.67735393583911
$row9[0,17] = 348.096185422552
$row9[0,18] = 0.501518814694381
$row9[0,19] = 0.501518814694381
$ws.Range("A9:T9").Value = $row9

$row10 = New-Object "object[,]" 1,20
$row10[0,0] = "MuSCs"
$row10[0,1] = "Clcf1"
$row10[0,2] = "Cntfr"
$row10[0,3] = "MuSCs"
$row10[0,4] = 3
$row10[0,5] = 1
$row10[0,6] = 6.825289333333334
$row10[0,7] = 20.475868
$row10[0,8] = 0.5251361975000832
$row10[0,9] = 0.5251361975000832
$row10[0,10] = 3
$row10[0,11] = 1
$row10[0,12] = 0.2284146666666667
$row10[0,13] = 0.685244
$row10[0,14] = 0.03849493351118214
$row10[0,15] = 0.03849493351118213
$row10[0,16] = 1.558996187976889
$row10[0,17] = 14.030965691792
$row10[0,18] = 0.02021508300708072
$row10[0,19] = 0.02021508300708071
$ws.Range("A10:T10").Value = $row10

$row11 = New-Object "object[,]" 1,20
$row11[0,0] = "Resolving-Mac"
$row11[0,1] = "Clcf1"
$row11[0,2] = "Cntfr"
$row11[0,3] = "ECs"
$row11[0,4] = 3
$row11[0,5] = 1
$row11[0,6] = 1.099475
$row11[0,7] = 3.298425
$row11[0,8] = 0.08459335458888541
$row11[0,9] = 0.08459335458888539
$row11[0,10] = 1
$row11[0,11] = 0.3333333333333333
$row11[0,12] = 0.03844333333333334
$row11[0,13] = 0.11533
$row11[0,14] = 0.006478890266597937
$row11[0,15] = 0.006478890266597936
$row11[0,16] = 0.04226748391666667
$row11[0,17] = 0.38040735525
$row11[0,18] = 0.0005480710616647975
$row11[0,19] = 0.0005480710616647974
$ws.Range("A11:T11").Value = $row11

$row12 = New-Object "object[,]" 1,20
$row12[0,0] = "Resolving-Mac"
$row12[0,1] = "Clcf1"
$row12[0,2] = "Cntfr"
$row12[0,3] = "FAPs"
$row12[0,4] = 3
$row12[0,5] = 1
$row12[0,6] = 1.099475
$row12[0,7] = 3.298425
$row12[0,8] = 0.08459335458888541
$row12[0,9] = 0.08459335458888539
$row12[0,10] = 3
$row12[0,11] = 1
$row12[0,12] = 5.666771333333333
$row12[0,13] = 17.000314
$row12[0,14] = 0.95502617622222
$row12[0,15] = 0.9550261762222199
$row12[0,16] = 6.230473411716666
$row12[0,17] = 56.07426070544999
$row12[0,18] = 0.08078886796683361
$row12[0,19] = 0.0807888679668336
$ws.Range("A12:T12").Value = $row12

$row13 = New-Object "object[,]" 1,20
$row13[0,0] = "Resolving-Mac"
$row13[0,1] = "Clcf1"
$row13[0,2] = "Cntfr"
$row13[0,3] = "MuSCs"
$row13[0,4] = 3
$row13[0,5] = 1
$row13[0,6] = 1.099475
$row13[0,7] = 3.298425
$row13[0,8] = 0.08459335458888541
$row13[0,9] = 0.08459335458888539
$row13[0,10] = 3
$row13[0,11] = 1
$row13[0,12] = 0.2284146666666667
$row13[0,13] = 0.685244
$row13[0,14] = 0.03849493351118214
$row13[0,15] = 0.03849493351118213
$row13[0,16] = 0.2511362156333333
$row13[0,17] = 2.2602259407
$row13[0,18] = 0.003256415560386998
$row13[0,19] = 0.003256415560386997
$ws.Range("A13:T13").Value = $row13
